$d = $word.ActiveDocument

$replacements = @(
    @("67×55=", "67×79="),
    @("68×40=", "98×24="),
    @("56×62=", "30×51="),
    @("47×25=", "66×30="),
    @("57×89=", "84×91="),
    @("46×30=", "53×81="),
    @("30×64=", "32×33="),
    @("34×38=", "60×43="),
    @("38×82=", "63×69="),
    @("57×95=", "74×38="),
    @("16×87=", "79×57="),
    @("98×56=", "26×13="),
    @("93×52=", "46×38="),
    @("88×57=", "93×18="),
    @("51×12=", "19×58="),
    @("91×89=", "59×95="),
    @("49×23=", "31×92="),
    @("42×33=", "64×25="),
    @("64×33=", "53×29="),
    @("34×59=", "54×63="),
    @("74×87=", "95×47="),
    @("53×56=", "94×24="),
    @("98×29=", "17×96="),
    @("11×11=", "88×85="),
    @("19×34=", "58×39=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
